# Apply the "Updated cryptos list" data refresh to Sheet1 (Coin / Link / Price / Volume(1h)).
# Row 13 and 14 swap content (Polkadot <-> WrappedEther) in addition to their Price/Volume updates.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell as plain text, even when it looks numeric
# (e.g. "1.000", "0.7690") so Excel does not silently coerce it to a Double and
# drop the formatting-significant trailing/leading characters. A leading apostrophe
# forces text entry; ClearFormats() then drops the "quote prefix" cell style Excel
# adds for that so the cell keeps its original (default) style index.
function Set-TextValue($cell, $text) {
    $ws.Range($cell).Formula = "'" + $text
    $ws.Range($cell).ClearFormats()
}

$ws.Range("D2").Value = "27.225.10"
$ws.Range("E2").Value = "  +0.86%  "
$ws.Range("D3").Value = "1.854.98"
$ws.Range("E3").Value = "  +1.51%  "
Set-TextValue "D4" "1.000"
$ws.Range("E4").Value = "  -0.48%  "
Set-TextValue "D5" "313.71"
Set-TextValue "D6" "1.000"
$ws.Range("E7").Value = "  +0.49%  "
Set-TextValue "D8" "0.3708"
$ws.Range("E8").Value = "  +0.27%  "
Set-TextValue "D9" "0.07297"
$ws.Range("E9").Value = "  -0.63%  "
Set-TextValue "D10" "0.8921"
$ws.Range("E10").Value = "  +1.87%  "
Set-TextValue "D12" "0.07866"
$ws.Range("E12").Value = "  -0.45%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D13" "5.401"
$ws.Range("E13").Value = "  +1.21%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.793.22"
$ws.Range("E14").Value = "  -1.33%  "
Set-TextValue "D15" "6.519"
$ws.Range("E15").Value = "  -0.60%  "
Set-TextValue "D16" "91.34"
$ws.Range("E16").Value = "  -0.17%  "
$ws.Range("E17").Value = "  -0.38%  "
Set-TextValue "D18" "0.000008916"
$ws.Range("E18").Value = "  +0.92%  "
$ws.Range("E19").Value = "  -0.31%  "
Set-TextValue "D20" "14.74"
$ws.Range("E20").Value = "  -0.31%  "
$ws.Range("D21").Value = "27.241.50"
$ws.Range("E21").Value = "  +0.85%  "
Set-TextValue "D22" "5.086"
$ws.Range("E22").Value = "  -0.29%  "
Set-TextValue "D23" "10.52"
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("D24").Value = "2.051.89"
$ws.Range("E24").Value = "  -1.23%  "
Set-TextValue "D25" "2.042"
$ws.Range("E25").Value = "  +9.76%  "
Set-TextValue "D26" "151.52"
$ws.Range("E26").Value = "  -0.46%  "
Set-TextValue "D27" "18.46"
$ws.Range("E27").Value = "  +0.17%  "
Set-TextValue "D28" "2.044"
$ws.Range("E28").Value = "  +0.25%  "
Set-TextValue "D29" "115.83"
$ws.Range("E29").Value = "  +0.11%  "
Set-TextValue "D30" "5.034"
$ws.Range("E30").Value = "  -1.41%  "
Set-TextValue "D31" "0.08832"
Set-TextValue "D32" "3.141"
$ws.Range("E32").Value = "  +6.09%  "
Set-TextValue "D33" "0.7690"
$ws.Range("E33").Value = "  +5.01%  "
$ws.Range("E34").Value = "  +3.07%  "
Set-TextValue "D35" "4.517"
$ws.Range("E35").Value = "  +1.72%  "
Set-TextValue "D36" "2.694"
$ws.Range("E36").Value = "  +9.15%  "
Set-TextValue "D37" "1.108"
$ws.Range("E37").Value = "  +3.11%  "
Set-TextValue "D38" "0.01941"
$ws.Range("E38").Value = "  -0.08%  "
Set-TextValue "D39" "0.05217"
$ws.Range("E39").Value = "  -0.08%  "
Set-TextValue "D40" "2.945"
$ws.Range("E40").Value = "  -0.51%  "
Set-TextValue "D41" "7.047"
$ws.Range("E41").Value = "  -0.70%  "
Set-TextValue "D42" "0.5106"
$ws.Range("E42").Value = "  -0.98%  "
$ws.Range("E43").Value = "  -0.01%  "
Set-TextValue "D44" "8.500"
$ws.Range("E44").Value = "  +4.22%  "
Set-TextValue "D45" "0.4797"
$ws.Range("E45").Value = "  -0.63%  "
Set-TextValue "D46" "10.35"
$ws.Range("E46").Value = "  +1.44%  "
Set-TextValue "D47" "1.000"
Set-TextValue "D48" "102.60"
$ws.Range("E48").Value = "  +0.73%  "
$ws.Range("E49").Value = "  +1.24%  "
Set-TextValue "D50" "0.06197"
$ws.Range("E50").Value = "  -0.15%  "
Set-TextValue "D51" "65.35"
$ws.Range("E51").Value = "  +0.88%  "
